# Error Calculations and Plots
# Applies the edits described in the commit diff to Sheet1:
#  - Update several "D" column (column E in the sheet) values, some
#    becoming blank (simulating missing data) and some becoming
#    populated (imputed values).
#  - Remove the two rows for "RM 232" and "SC 92" entirely.
#  - Fill in the "D" value for "SC 193".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
Write-Host ("Initial last row: " + $lastRow)

# Helper: find the row number whose column A equals a given id text.
function Find-RowById($ws, $id, $maxRow) {
    for ($r = 2; $r -le $maxRow; $r++) {
        $cellVal = $ws.Cells.Item($r, 1).Value()
        if ($cellVal -eq $id) {
            return $r
        }
    }
    return -1
}

# --- Column E (4th data column) value changes, matched by row ID ---
$rowRM8  = Find-RowById $ws "RM 8" $lastRow
$rowRM14 = Find-RowById $ws "RM 14" $lastRow
$rowRM135 = Find-RowById $ws "RM 135" $lastRow
$rowRM140 = Find-RowById $ws "RM 140" $lastRow
$rowSC193 = Find-RowById $ws "SC 193" $lastRow

Write-Host ("RM 8 row: " + $rowRM8)
Write-Host ("RM 14 row: " + $rowRM14)
Write-Host ("RM 135 row: " + $rowRM135)
Write-Host ("RM 140 row: " + $rowRM140)
Write-Host ("SC 193 row: " + $rowSC193)

# RM 8: E was blank -> now -5.7
$ws.Cells.Item($rowRM8, 5).Value = -5.7

# RM 14: E was -5 -> now blank
$ws.Cells.Item($rowRM14, 5).ClearContents()

# RM 135: E was blank -> now -8.699999999999999
$ws.Cells.Item($rowRM135, 5).Value = -8.699999999999999

# RM 140: E was -7 -> now blank
$ws.Cells.Item($rowRM140, 5).ClearContents()

# SC 193: E was blank -> now -6.4
$ws.Cells.Item($rowSC193, 5).Value = -6.4

# --- Remove the rows for "RM 232" and "SC 92" entirely ---
$rowRM232 = Find-RowById $ws "RM 232" $lastRow
$rowSC92  = Find-RowById $ws "SC 92" $lastRow
Write-Host ("RM 232 row: " + $rowRM232)
Write-Host ("SC 92 row: " + $rowSC92)

# Delete the higher-numbered row first so the other row index stays valid.
if ($rowRM232 -gt $rowSC92) {
    $ws.Rows($rowRM232).Delete()
    $ws.Rows($rowSC92).Delete()
} else {
    $ws.Rows($rowSC92).Delete()
    $ws.Rows($rowRM232).Delete()
}

$finalLastRow = $ws.UsedRange.Rows.Count
Write-Host ("Final last row: " + $finalLastRow)
